# "fixes to summary formulas"
# Expands the summary block (rows 1-6) to account for two more fee
# columns (International Fees -> col L, Other Fee -> col M), widens the
# affected SUM ranges out to the full sheet, and re-labels / reflows the
# Total Expenses / Total Profit / ROI rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 --------------------------------------------------------
$ws.Range("H1").Formula = "=SUM(J8:J1048576)"

# ---- Row 2 ----------------------------------------------------------
$ws.Range("F2").Formula = "=SUM(F8:F1048576)"
$ws.Range("H2").Formula = "=SUM(K8:K1048576)"
$ws.Range("J2").Formula = "=F4+H1+H2+H3+H4+H5"

# ---- Row 3 ----------------------------------------------------------
$ws.Range("F3").Formula = "=SUM(H8:H1048576)"
$ws.Range("G3").Value = "International Fees"
$ws.Range("H3").Formula = "=SUM(L8:L1048576)"
$ws.Range("I3").Value = "Total Expenses"
$ws.Range("J3").NumberFormat = "_(""$""* #,##0.00_);_(""$""* (#,##0.00);_(""$""* ""-""??_);_(@_)"
$ws.Range("J3").Value = 0
$ws.Range("K3").NumberFormat = "_(""$""* #,##0.00_);_(""$""* (#,##0.00);_(""$""* ""-""??_);_(@_)"
$ws.Range("L3").Font.Bold = $true
$ws.Range("M3").Font.Bold = $true
$ws.Range("M3").HorizontalAlignment = -4152

# ---- Row 4 ----------------------------------------------------------
$ws.Range("F4").Formula = "=SUM(E8:E1048576)"
$ws.Range("G4").Value = "Other Fee"
$ws.Range("G4").Font.Bold = $true
$ws.Range("H4").Formula = "=SUM(M8:M1048576)"
$ws.Range("I4").Value = "Total Profit"
$ws.Range("J4").Formula = "=J1+J2+J3"
$ws.Range("K4").NumberFormat = "_(""$""* #,##0.00_);_(""$""* (#,##0.00);_(""$""* ""-""??_);_(@_)"
$ws.Range("L4").NumberFormat = "_(""$""* #,##0.00_);_(""$""* (#,##0.00);_(""$""* ""-""??_);_(@_)"
$ws.Range("M4").NumberFormat = "_(""$""* #,##0.00_);_(""$""* (#,##0.00);_(""$""* ""-""??_);_(@_)"

# ---- Row 5 ----------------------------------------------------------
$ws.Range("C5").NumberFormat = "@"
$ws.Range("G5").Value = "S&H Cost"
$ws.Range("H5").Formula = "=SUM(I8:I1048576)"
$ws.Range("I5").Value = "ROI"
$ws.Range("J5").Formula = "=IFERROR(J4/J1,"" "")"

# ---- Column widths / formatting for col G ---------------------------
$ws.Columns("G").ColumnWidth = 19.7109375

# ---- Selection / view ------------------------------------------------
$ws.Range("F17").Select()
